$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.091.31"
$ws.Range("E2").Value = "  -2.32%  "

$ws.Range("D3").Value = "2.491.85"
$ws.Range("E3").Value = "  -4.37%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "553.75"
$ws.Range("E5").Value = "  -3.34%  "

$ws.Range("D6").Value = "147.25"
$ws.Range("E6").Value = "  -4.89%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  -2.67%  "

$ws.Range("D9").Value = "2.488.59"
$ws.Range("E9").Value = "  -4.46%  "

$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  -7.74%  "

$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").Value = "5.42"
$ws.Range("E12").Value = "  -6.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.360"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.67%  "

$ws.Range("D14").Value = "26.27"
$ws.Range("E14").Value = "  -7.02%  "

$ws.Range("D15").Value = "2.941.01"
$ws.Range("E15").Value = "  -4.34%  "

$ws.Range("D16").Value = "0.0000165"
$ws.Range("E16").Value = "  -6.75%  "

$ws.Range("D17").Value = "61.951.38"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").Value = "2.497.14"
$ws.Range("E18").Value = "  -4.10%  "

$ws.Range("D19").Value = "11.27"
$ws.Range("E19").Value = "  -5.73%  "

$ws.Range("D20").Value = "7.03"
$ws.Range("E20").Value = "  -6.03%  "

$ws.Range("D21").Value = "4.24"
$ws.Range("E21").Value = "  -6.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "324.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.46%  "

$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "64.68"
$ws.Range("E24").Value = "  -3.50%  "

$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").Value = "0.0000103"
$ws.Range("E26").Value = "  -4.05%  "

$ws.Range("D27").Value = "2.640.72"
$ws.Range("E27").Value = "  -3.28%  "

$ws.Range("D28").Value = "1.51"
$ws.Range("E28").Value = "  -3.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -6.79%  "

$ws.Range("D30").Value = "541.51"
$ws.Range("E30").Value = "  -7.88%  "

$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").Value = "7.73"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").Value = "0.151"
$ws.Range("E33").Value = "  -5.62%  "

$ws.Range("D34").Value = "1.91"
$ws.Range("E34").Value = "  -6.97%  "

$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  -7.65%  "

$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  -9.03%  "

$ws.Range("D37").Value = "4.88"
$ws.Range("E37").Value = "  -9.11%  "

$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.380"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.90%  "

$ws.Range("D40").Value = "18.67"
$ws.Range("E40").Value = "  -4.90%  "

$ws.Range("D41").Value = "147.12"
$ws.Range("E41").Value = "  -5.12%  "

$ws.Range("D42").Value = "1.71"
$ws.Range("E42").Value = "  -7.73%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "40.82"
$ws.Range("E44").Value = "  -1.49%  "

$ws.Range("D45").Value = "2.34"
$ws.Range("E45").Value = "  -3.93%  "

$ws.Range("D46").Value = "149.11"
$ws.Range("E46").Value = "  -4.32%  "

$ws.Range("D47").Value = "3.64"
$ws.Range("E47").Value = "  -6.15%  "

$ws.Range("D48").Value = "21.52"
$ws.Range("E48").Value = "  -7.05%  "

$ws.Range("D49").Value = "0.0543"
$ws.Range("E49").Value = "  -7.28%  "

$ws.Range("D50").Value = "0.594"
$ws.Range("E50").Value = "  -5.38%  "

$ws.Range("D51").Value = "0.0955"
$ws.Range("E51").Value = "  -4.70%  "
